# Update the "想去人数" (interested-count) figures that were refreshed in
# this data snapshot. Same underlying events appear on sheet "展览"
# (F column), are mirrored on "全部类型", and a couple of matching updates
# also land on "演出" and "本地生活".

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsLocal   = $wb.Worksheets.Item("本地生活")
$wsAll     = $wb.Worksheets.Item("全部类型")

# --- 展览 (Exhibition) sheet ---
$wsExhibit.Range("F4").Value  = 423
$wsExhibit.Range("F5").Value  = 8560
$wsExhibit.Range("F7").Value  = 10742
$wsExhibit.Range("F22").Value = 1821
$wsExhibit.Range("F23").Value = 80
$wsExhibit.Range("F24").Value = 565
$wsExhibit.Range("F27").Value = 67
$wsExhibit.Range("F28").Value = 585
$wsExhibit.Range("F30").Value = 1194
$wsExhibit.Range("F33").Value = 1417
$wsExhibit.Range("F34").Value = 446
$wsExhibit.Range("F35").Value = 345
$wsExhibit.Range("F36").Value = 289
$wsExhibit.Range("F37").Value = 25
$wsExhibit.Range("F39").Value = 517
$wsExhibit.Range("F40").Value = 349
$wsExhibit.Range("F41").Value = 99
$wsExhibit.Range("F43").Value = 643
$wsExhibit.Range("F45").Value = 114
$wsExhibit.Range("F46").Value = 104

# --- 演出 (Show) sheet ---
$wsShow.Range("F3").Value = 28
$wsShow.Range("F6").Value = 47

# --- 本地生活 (Local life) sheet ---
$wsLocal.Range("F3").Value = 2805

# --- 全部类型 (All types) sheet, mirrors the rows above ---
$wsAll.Range("F7").Value  = 28
$wsAll.Range("F8").Value  = 423
$wsAll.Range("F9").Value  = 8560
$wsAll.Range("F11").Value = 10742
$wsAll.Range("F19").Value = 1821
$wsAll.Range("F20").Value = 80
$wsAll.Range("F21").Value = 565
$wsAll.Range("F24").Value = 67
$wsAll.Range("F26").Value = 585
$wsAll.Range("F28").Value = 1194
$wsAll.Range("F34").Value = 1417
$wsAll.Range("F35").Value = 446
$wsAll.Range("F37").Value = 345
$wsAll.Range("F39").Value = 517
$wsAll.Range("F40").Value = 349
$wsAll.Range("F41").Value = 99
$wsAll.Range("F46").Value = 643
$wsAll.Range("F48").Value = 114
$wsAll.Range("F49").Value = 105

$wb.Save()
